# Update forest data - 2026-01-29 12:29
#
# Moves the two existing rows from the "New" sheet into the
# "Previously added" sheet (appended at the bottom), then populates the
# "New" sheet with six freshly scraped listings.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# 1. Capture the hyperlink targets of the two rows currently on "New"
#    before we touch anything (row 2 -> A440, row 3 -> A441).
# ---------------------------------------------------------------------
$oldLink2 = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/jersikas-pag/fpncc.html"
$oldLink3 = "https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/zvardes-pag/kjfkn.html"

# ---------------------------------------------------------------------
# 2. Copy the two existing "New" rows (values + formatting) down to the
#    bottom of "Previously added" as rows 440 and 441.
# ---------------------------------------------------------------------
$ws2.Range("A2:F2").Copy($ws1.Range("A440:F440"))
$ws2.Range("A3:F3").Copy($ws1.Range("A441:F441"))

$ws1.Hyperlinks.Add($ws1.Range("A440"), $oldLink2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A441"), $oldLink3) | Out-Null

# Restore the original hyperlink cell formatting (Hyperlinks.Add re-styles
# the anchor cell with a fresh auto-generated style) by re-pasting the
# formats from the row directly above, which already carries the correct
# hyperlink look.
$ws1.Range("A439:F439").Copy()
$ws1.Range("A440:F441").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 3. Clear out the old rows 2:3 on "New" (data + hyperlinks) so the sheet
#    can be repopulated with the newly scraped listings.
# ---------------------------------------------------------------------
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Range("A2:F3").ClearContents()

# Re-apply the standard row formatting (link / price / region / area /
# cadastre-code / date) across the six new rows by tiling the format of
# the template row from "Previously added".
$ws1.Range("A439:F439").Copy()
$ws2.Range("A2:F7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 4. Fill in the six new listings.
# ---------------------------------------------------------------------

# Row 2
$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/lazdukalna-pag/gboid.html"
$ws2.Range("B2").Value = "15 000 €"
$ws2.Range("C2").Value = "Balvi un raj."
$ws2.Range("D2").Value = "2 ha."
$ws2.Range("E2").Value = "'38640090004"
$ws2.Range("F2").Value = 46050.70763888889
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/lazdukalna-pag/gboid.html") | Out-Null

# Row 3
$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/mbgnx.html"
$ws2.Range("B3").Value = "30 000 €"
$ws2.Range("C3").Value = "Bauska un raj."
$ws2.Range("D3").Value = "10 ha."
$ws2.Range("E3").Value = "'40640110604"
$ws2.Range("F3").Value = 46051.54236111111
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/mbgnx.html") | Out-Null

# Row 4
$ws2.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/salienas-pag/ojhkp.html"
$ws2.Range("B4").Value = "12 800 €"
$ws2.Range("C4").Value = "Daugavpils un raj."
$ws2.Range("D4").Value = "1 ha."
$ws2.Range("E4").Value = "'44840040659"
$ws2.Range("F4").Value = 46050.72222222222
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/salienas-pag/ojhkp.html") | Out-Null

# Row 5
$ws2.Range("A5").Value = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/belavas-pag/bxxehb.html"
$ws2.Range("B5").Value = "26 000 €"
$ws2.Range("C5").Value = "Gulbene un raj."
$ws2.Range("D5").Value = "3 ha."
$ws2.Range("E5").Value = "'50440070082"
$ws2.Range("F5").Value = 46050.71736111111
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/belavas-pag/bxxehb.html") | Out-Null

# Row 6
$ws2.Range("A6").Value = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/belavas-pag/injkh.html"
$ws2.Range("B6").Value = "140 000 €"
$ws2.Range("C6").Value = "Gulbene un raj."
$ws2.Range("D6").Value = "43 ha."
$ws2.Range("E6").Value = "'50440060014"
$ws2.Range("F6").Value = 46050.70972222222
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/belavas-pag/injkh.html") | Out-Null

# Row 7
$ws2.Range("A7").Value = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rusonas-pag/ixeng.html"
$ws2.Range("B7").Value = "211 000 €"
$ws2.Range("C7").Value = "Preiļi un raj."
$ws2.Range("D7").Value = "91 ha."
$ws2.Range("E7").Value = "7670 011 0032 011 00"
$ws2.Range("F7").Value = 46051.17083333334
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rusonas-pag/ixeng.html") | Out-Null

# ---------------------------------------------------------------------
# 5. Re-paint formatting once more over the whole new-row block. This
#    scrubs the transient "quote prefix" cell style that Excel applies
#    when text looking like a number is entered with a leading apostrophe
#    (columns E2:E6), restoring the plain template style (s=4) while
#    keeping the cell's string value/type intact.
# ---------------------------------------------------------------------
$ws1.Range("A439:F439").Copy()
$ws2.Range("A2:F7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 6. Re-apply the date/time number format on column F (PasteSpecial of
#    formats already set numFmt 164, but make sure explicitly).
# ---------------------------------------------------------------------
$ws1.Range("F440:F441").NumberFormat = "dd.mm.yyyy hh:mm"
$ws2.Range("F2:F7").NumberFormat = "dd.mm.yyyy hh:mm"

Write-Host "Done: Previously added rows = $($ws1.UsedRange.Rows.Count); New rows = $($ws2.UsedRange.Rows.Count)"
